# Update the USD Amount value in T2 and move the active selection to T3,
# matching the authored edit captured in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T2").Value = 163512
$ws.Range("T3").Select()
